$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

# Row 2
$ws.Range("A2").Value = 111749343
$ws.Range("B2").Value = 78107
$ws.Range("D2").Value = "NT"
$ws.Range("E2").Value = 6453
$ws.Range("F2").Value = "Vedskivlav"
$ws.Range("G2").Value = "Hertelidea botryosa"
$ws.Range("H2").Value = "(Fr.) Printzen & Kantvilas"
$ws.Range("Q2").Value = 575415.2450877089
$ws.Range("R2").Value = 6509807.674603676
$ws.Range("S2").Value = 1

# Row 3
$ws.Range("A3").Value = 111747705
$ws.Range("B3").Value = 93067
$ws.Range("E3").Value = 2810
$ws.Range("F3").Value = "Västlig hakmossa"
$ws.Range("G3").Value = "Rhytidiadelphus loreus"
$ws.Range("H3").Value = "(Hedw.) Warnst."
$ws.Range("Q3").Value = 575459.4222356658
$ws.Range("R3").Value = 6509864.113963567
$ws.Range("S3").Value = 2

# Row 4
$ws.Range("A4").Value = 111747186
$ws.Range("P4").Value = "Lilla gruvan (Lilla gruvan), Ög"
$ws.Range("Q4").Value = 575435.6246570286
$ws.Range("R4").Value = 6509856.898648335
$ws.Range("S4").Value = 2

# Row 5
$ws.Range("A5").Value = 111749883
$ws.Range("P5").Value = "Älgsjöhåll (Älgsjöhåll), Ög"
$ws.Range("Q5").Value = 575336.5075504743
$ws.Range("R5").Value = 6509789.003789719
$ws.Range("S5").Value = 1

# Row 7
$ws.Range("A7").Value = 111749097
$ws.Range("B7").Value = 93388
$ws.Range("D7").Value = "LC"
$ws.Range("E7").Value = 2180
$ws.Range("F7").Value = "Blåmossa"
$ws.Range("G7").Value = "Leucobryum glaucum"
$ws.Range("H7").Value = "(Hedw.) Ångstr."
$ws.Range("P7").Value = "Lilla gruvan (Lilla gruvan), Ög"
$ws.Range("Q7").Value = 575501.7342092508
$ws.Range("R7").Value = 6509775.591426332
$ws.Range("S7").Value = 3

# Row 8
$ws.Range("A8").Value = 111749860
$ws.Range("P8").Value = "Älgsjöhåll (Älgsjöhåll), Ög"
$ws.Range("Q8").Value = 575356.6078101217
$ws.Range("R8").Value = 6509772.251964441

# Row 9
$ws.Range("A9").Value = 111749897
$ws.Range("Q9").Value = 575336.6687912485
$ws.Range("R9").Value = 6509780.695668718
